$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old rows 7-15 held extra leads that are no longer needed; remove them
# so the sheet only keeps the header plus 5 data rows (A1:G6).
$ws.Range("A7:G15").EntireRow.Delete()

# Replace the remaining lead rows (2-6) with the refreshed dataset (new
# names/emails/proxy as a failsafe for the email body on smtp).
$data = @(
    @("Franklin","Logans","frankinsonloslogansi@gmail.com","kdejqyra","185.24.233.182:4006","8GbKtEpRUr29jbg5","TMwprA4NyqSKxc6V"),
    @("Frederick","Salamon","fredericksongatsalamon@gmail.com","kyqcxzun","185.24.233.182:4007","8GbKtEpRUr29jbg5","TMwprA4NyqSKxc6V"),
    @("Franklin","Morrison","frankinsonleemorison443@gmail.com","qzbsvkrp","185.24.233.182:4008","8GbKtEpRUr29jbg5","TMwprA4NyqSKxc6V"),
    @("Frederick","Vartinson","fredriksonkarlosvartinson@gmail.com","enpxdtwh","185.24.233.182:4009","8GbKtEpRUr29jbg5","TMwprA4NyqSKxc6V"),
    @("Franklin","Johnsson","frankinsonleejhonson5@gmail.com","vznrcpwb","185.24.233.182:4006","8GbKtEpRUr29jbg5","TMwprA4NyqSKxc6V")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($row, $j + 1).Value = $rowData[$j]
    }
}

# Match the author's final selection/active cell in the saved workbook.
$ws.Range("I12").Select() | Out-Null
